# Generate Report for Handoff
# Regenerate the handoff report: bump Priority from "low" to "ht" for the
# files that are "Ready for handoff", and refresh the handoff timestamps.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# zh-cn: rows 4-7 are the "Ready for handoff" files. Priority (col E) flips
# from low to ht, and the Latest Handoff Datetime (col H) is refreshed.
foreach ($r in 4..7) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-22 00:43:29"
}

# de-de: same four rows get the priority bump; their handoff datetime
# (col H) shares its text with the Overview's "Latest HO Xliff Generate
# Date" column and is refreshed to a later timestamp.
foreach ($r in 4..7) {
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-22 00:43:33"
}

# Overview: rows 4-7, col G ("Latest HO Xliff Generate Date") mirrors the
# de-de handoff datetime above.
foreach ($r in 4..7) {
    $overview.Range("G$r").Value = "2016-08-22 00:43:33"
}
